# Update "想去人数" (F column) counts across the four worksheets to the
# values captured in the newer data pull.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 295
$ws.Range("F3").Value = 495
$ws.Range("F4").Value = 37
$ws.Range("F7").Value = 1227
$ws.Range("F8").Value = 382
$ws.Range("F10").Value = 350
$ws.Range("F11").Value = 8215
$ws.Range("F13").Value = 9885
$ws.Range("F26").Value = 389
$ws.Range("F27").Value = 1696
$ws.Range("F28").Value = 41
$ws.Range("F29").Value = 491
$ws.Range("F30").Value = 315
$ws.Range("F32").Value = 47
$ws.Range("F33").Value = 551
$ws.Range("F34").Value = 996
$ws.Range("F35").Value = 12
$ws.Range("F39").Value = 325
$ws.Range("F41").Value = 120
$ws.Range("F42").Value = 487
$ws.Range("F43").Value = 310
$ws.Range("F44").Value = 62
$ws.Range("F46").Value = 106
$ws.Range("F48").Value = 21
$ws.Range("F49").Value = 23

# --- Sheet "演出" ---------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 33

# --- Sheet "本地生活" -----------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2765
$ws.Range("F5").Value = 192

# --- Sheet "全部类型" -----------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 295
$ws.Range("F3").Value = 495
$ws.Range("F4").Value = 33
$ws.Range("F6").Value = 192
$ws.Range("F7").Value = 37
$ws.Range("F10").Value = 1227
$ws.Range("F11").Value = 382
$ws.Range("F16").Value = 8215
$ws.Range("F18").Value = 9885
$ws.Range("F24").Value = 1696
$ws.Range("F25").Value = 41
$ws.Range("F26").Value = 315
$ws.Range("F28").Value = 47
$ws.Range("F30").Value = 551
$ws.Range("F40").Value = 120
$ws.Range("F41").Value = 487
$ws.Range("F43").Value = 62
$ws.Range("F48").Value = 21
$ws.Range("F49").Value = 23
